$wb = $excel.ActiveWorkbook

# Rename SOFI_PROP_* -> SAND_PROP_* across the whole workbook (all sheets).
# (This is the "renamed all SOFI to SAND" commit - it touches every sheet
# header row that shows the SOFI_PROP_BASE/MDB/RANGE/SAVE/ACCESS group
# labels.)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("SOFI", "SAND") | Out-Null
}

# Update the saved cursor position on a couple of sheets.
$wsDevice = $wb.Worksheets.Item(5)
$wsDevice.Range("Q17").Select() | Out-Null

$wsAi = $wb.Worksheets.Item(9)
$wsAi.Range("N17").Select() | Out-Null

# "test" becomes the active/selected sheet (workbook activeTab points at it,
# and it gets tabSelected="1" while "template" loses it), with its
# bottom-right pane selection on O15.
$wsTest = $wb.Worksheets.Item(10)
$wsTest.Activate() | Out-Null
$wsTest.Range("O15").Select() | Out-Null
